# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "59.416.11"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "2.408.86"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  +7.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.86%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "2.822.91"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "59.234.84"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000140"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.35%  "
$ws.Range("D17").Value = "2.407.02"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.45%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.993"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.96%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "40.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "295.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0964"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.572"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0225"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.397"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
